# Auto-generated edit script applying the scheduled-runner price/profit refresh
# across the Table_* sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 142.78572
$ws.Range("I4").Value = 46.363636
$ws.Range("K4").Value = 46.363636
$ws.Range("M4").Value = 67.636364
$ws.Range("H74").Value = 6729.25
$ws.Range("I74").Value = 4042.7144
$ws.Range("K74").Value = 4042.7144
$ws.Range("M74").Value = -3106.7144
$ws.Range("H76").Value = 4391.76
$ws.Range("I76").Value = 3808.1052
$ws.Range("K76").Value = 3808.1052
$ws.Range("M76").Value = -3493.1052
$ws.Range("H77").Value = 6729.25
$ws.Range("I77").Value = 4042.7144
$ws.Range("K77").Value = 20213.572
$ws.Range("M77").Value = -15533.572
$ws.Range("H79").Value = 4391.76
$ws.Range("I79").Value = 3808.1052
$ws.Range("K79").Value = 3808.1052
$ws.Range("M79").Value = -2716.1052
$ws.Range("H88").Value = 5906.1333
$ws.Range("J88").Value = 5829.3
$ws.Range("L88").Value = 5829.3
$ws.Range("N88").Value = -6641.3
$ws.Range("H91").Value = 5906.1333
$ws.Range("J91").Value = 5829.3
$ws.Range("L91").Value = 5829.3
$ws.Range("N91").Value = -8637.299999999999
$ws.Range("H100").Value = 5670.3335
$ws.Range("I100").Value = 2594.3333
$ws.Range("J100").Value = 8746.333000000001
$ws.Range("K100").Value = 2594.3333
$ws.Range("L100").Value = 8746.333000000001
$ws.Range("M100").Value = -2053.3333
$ws.Range("N100").Value = -9828.333000000001
$ws.Range("H129").Value = 1485.579
$ws.Range("I129").Value = 1170.5
$ws.Range("J129").Value = 3166
$ws.Range("K129").Value = 3511.5
$ws.Range("L129").Value = 9498
$ws.Range("M129").Value = 1488.5
$ws.Range("N129").Value = -19498
$ws.Range("H132").Value = 2268.2144
$ws.Range("I132").Value = 2135.077
$ws.Range("K132").Value = 6405.231000000001
$ws.Range("M132").Value = -3875.231000000001
$ws.Range("H137").Value = 2742.35
$ws.Range("I137").Value = 1713.2
$ws.Range("K137").Value = 5139.6
$ws.Range("M137").Value = -2589.6

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3379.92
$ws.Range("I32").Value = 3386.0852
$ws.Range("K32").Value = 3386.0852
$ws.Range("M32").Value = -3099.0852
$ws.Range("H102").Value = 4702.467
$ws.Range("I102").Value = 4394.7
$ws.Range("K102").Value = 4394.7
$ws.Range("M102").Value = -2772.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1499.2916
$ws.Range("I105").Value = 1499.2916
$ws.Range("K105").Value = 1499.2916
$ws.Range("M105").Value = 247.7084
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4185.8
$ws.Range("I86").Value = 4309.6665
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 4309.6665
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -3186.6665
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 4185.8
$ws.Range("I89").Value = 4309.6665
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 21548.3325
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -15932.3325
$ws.Range("N89").Value = -31232

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 340
$ws.Range("I38").Value = 120
$ws.Range("K38").Value = 360
$ws.Range("M38").Value = -13
$ws.Range("H109").Value = 4214.263
$ws.Range("I109").Value = 2766.111
$ws.Range("J109").Value = 5517.6
$ws.Range("K109").Value = 8298.332999999999
$ws.Range("L109").Value = 16552.8
$ws.Range("M109").Value = -7258.332999999999
$ws.Range("N109").Value = -18632.8
$ws.Range("I117").Value = 1830.75
$ws.Range("J117").Value = 4415.357
$ws.Range("K117").Value = 5492.25
$ws.Range("L117").Value = 13246.071
$ws.Range("M117").Value = -2050.25
$ws.Range("N117").Value = -20130.071

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19199
$ws.Range("J33").Value = 19199
$ws.Range("L33").Value = 19199
$ws.Range("N33").Value = -19703
$ws.Range("H36").Value = 49999
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H70").Value = 4624.875
$ws.Range("I70").Value = 3999.6667
$ws.Range("K70").Value = 3999.6667
$ws.Range("M70").Value = -3729.6667
$ws.Range("H73").Value = 4624.875
$ws.Range("I73").Value = 3999.6667
$ws.Range("K73").Value = 3999.6667
$ws.Range("M73").Value = -3063.6667
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 1447.0233
$ws.Range("I132").Value = 1455.2972
$ws.Range("K132").Value = 4365.8916
$ws.Range("M132").Value = -1835.8916

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 85466.25
$ws.Range("I22").Value = 201640.2
$ws.Range("J22").Value = 2484.8572
$ws.Range("K22").Value = 201640.2
$ws.Range("L22").Value = 2484.8572
$ws.Range("M22").Value = -201345.2
$ws.Range("N22").Value = -3074.8572
$ws.Range("H27").Value = 85466.25
$ws.Range("I27").Value = 201640.2
$ws.Range("J27").Value = 2484.8572
$ws.Range("K27").Value = 201640.2
$ws.Range("L27").Value = 2484.8572
$ws.Range("M27").Value = -201533.2
$ws.Range("N27").Value = -2698.8572
$ws.Range("H61").Value = 2292.75
$ws.Range("I61").Value = 1548.4
$ws.Range("K61").Value = 1548.4
$ws.Range("M61").Value = -1346.4
$ws.Range("H93").Value = 2962.0557
$ws.Range("I93").Value = 3076.3333
$ws.Range("K93").Value = 3076.3333
$ws.Range("M93").Value = -1828.3333
$ws.Range("H113").Value = 2292.75
$ws.Range("I113").Value = 1548.4
$ws.Range("K113").Value = 1548.4
$ws.Range("M113").Value = 621.5999999999999
$ws.Range("H122").Value = 3502
$ws.Range("I122").Value = 3206.56
$ws.Range("J122").Value = 4557.143
$ws.Range("K122").Value = 9619.68
$ws.Range("L122").Value = 13671.429
$ws.Range("M122").Value = -7169.68
$ws.Range("N122").Value = -18571.429
$ws.Range("H136").Value = 5854.5
$ws.Range("I136").Value = 4725.4
$ws.Range("K136").Value = 14176.2
$ws.Range("M136").Value = -11626.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1675.7715
$ws.Range("J113").Value = 3525.182
$ws.Range("L113").Value = 10575.546
$ws.Range("N113").Value = -14915.546
$ws.Range("H126").Value = 2394
$ws.Range("I126").Value = 2299.5881
$ws.Range("K126").Value = 6898.7643
$ws.Range("M126").Value = -4428.7643
